# Updates cryptos list values (Price / Volume(1h) columns) and fixes the
# RenderToken/Bittensor row order, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.523.51"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.478.35"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'489.08"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'148.90"
$ws.Range("E6").Value = "  +6.79%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.484.88"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  +5.99%  "
$ws.Range("D11").Value = "'0.0979"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "2.912.75"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "56.457.14"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "'20.95"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "2.489.56"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").Value = "'10.15"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "'318.13"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Value = "'58.49"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'0.162"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "2.598.24"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "'7.60"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "0.0₃0785"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'149.13"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").Value = "'5.17"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "'0.869"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").Value = "'1.39"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("D40").Value = "'33.78"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "'0.995"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "'0.0554"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "'0.608"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'264.43"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.79"
$ws.Range("E46").Value = "  +8.71%  "
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "'10.19"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'17.64"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "1.890.23"
$ws.Range("E51").Value = "  -4.37%  "
